# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (E3) and "Correspond Handback DateTime" (H3)
# for the second file row (77c44467-...) on both the zh-cn and de-de report sheets,
# reflecting a newly generated handback report with fresh timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-20 14:43:29"
$wsZhCn.Range("H3").Value = "2016-03-20 14:43:55"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-20 14:43:32"
$wsDeDe.Range("H3").Value = "2016-03-20 14:44:01"
